$d = $word.ActiveDocument

# Paragraph 1 is the hidden merge-field marker paragraph:
#   "**ID__AFFARS_pgi_5304_topic_19__ID** "   (token run + trailing space run)
# It needs to become:
#   "**ID__AFFARS_SMC_PGI_5304_1601__ID**"    (token run only, no trailing space)
# plus a tighter left indent and a (line-less) paragraph border used only
# for its spacing attributes, matching paragraph 3 later in the document.
$p1 = $d.Paragraphs(1)

$d.Content.Find.Execute("**ID__AFFARS_pgi_5304_topic_19__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SMC_PGI_5304_1601__ID**", 2)

# Remove the now-orphaned trailing-space run at the end of the paragraph.
$tailRange = $p1.Range.Duplicate
$foundSpace = $tailRange.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundSpace) {
    $tailRange.Delete()
}

# Tighten the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Add paragraph-border spacing (no visible line) around the paragraph.
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromRight = 5

Write-Output "Updated paragraph 1 token, indent and border spacing."
